# "fix error in linux"
# The SqlServer.xlsx Ini sheet hard-coded a loopback IP (127.0.0.1) for the
# SQL IP column, which only worked when client and server were on the same
# host. Replace it with a real LAN address (192.168.1.113) so the config
# also works for a separate Linux box, mark the cell as Text (so Excel
# doesn't try to reinterpret the dotted value), widen column C to fit the
# longer string, and leave the selection on the cell that was edited.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C2 holds the SqlIP value - swap the loopback address for the real one.
$ws.Range("C2").Value = "192.168.1.113"

# Keep it as text (same numeric format as the neighbouring ServerID cell)
# so Excel doesn't mangle the dotted-quad value.
$ws.Range("C2").NumberFormat = "@"

# Column C used to share column B's width (10.5); now that it carries a
# longer IP string it needs its own, wider column (stored width 15).
$ws.Columns.Item(3).ColumnWidth = 100/7

# Leave the selection on the cell that was just edited.
$ws.Range("C2").Select() | Out-Null
